$d = $word.ActiveDocument

# 1. Update the execution date/location sentence
$d.Content.Find.Execute(
    "THIS DEED OF PARTNERSHIP is executed on this 01/01/2024 at Ratnagiri by and between:", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "THIS DEED OF PARTNERSHIP is executed on this 2025-05-15 at Ratnagiri, Maharashtra by and between:", `
    2)

# 2. Partner No. 1 name -> [Full name]
$d.Content.Find.Execute(
    "1. Advait Milind Kulkarni, Son of Milind Shashikant Kulkarni, Age 25, residing at 557/H1,Thiba Palace Road ,AnandNagar, Ratnagiri, Maharashtra (Hereinafter referred to as Partner No. 1)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "1. [Full name], Son of Milind Shashikant Kulkarni, Age 25, residing at 557/H1,Thiba Palace Road ,AnandNagar, Ratnagiri, Maharashtra (Hereinafter referred to as Partner No. 1)", `
    2)

# 3. Partner No. 2 name -> [Full name] and new address
$d.Content.Find.Execute(
    "2. Tanmay Abhay Joshi, Son of Abhay Joshi, Age 25, residing at 557/H1,Thiba Palace Road ,AnandNagar, Ratnagiri, Maharashtra (Hereinafter referred to as Partner No. 2)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "2. [Full name], Son of Abhay Joshi, Age 25, residing at Omkar Sanjiwani Apartment Ratnagiri 415612 (Hereinafter referred to as Partner No. 2)", `
    2)

# 4. Commencement date
$d.Content.Find.Execute(
    "The partnership shall commence on 01/01/2025 and shall be a Partnership at Will.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "The partnership shall commence on 2025-05-30 and shall be a Partnership at Will.", `
    2)

# 5. Replace severability clause with governing law / jurisdiction clause
$d.Content.Find.Execute(
    "Severability: If any provision of this agreement is held to be invalid or unenforceable by a court of competent jurisdiction, the remaining provisions shall continue to be valid and enforceable to the fullest extent permitted by law.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "This agreement shall be governed by and construed in accordance with the laws of the State of [insert state], without regard to its conflict of laws principles. Any disputes arising out of or related to this agreement shall be exclusively resolved in the state and federal courts located in [insert city, state], and the parties hereby consent to the personal jurisdiction of such courts.", `
    2)

# 6 & 7. Signature table cells: full names -> [Full name]
$table = $d.Tables.Item(1)
$table.Cell(2, 1).Range.Text = "[Full name]"
$table.Cell(2, 2).Range.Text = "[Full name]"
